$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: new "ford" entry (same text as row 2) with a red font, plus
#     "new model" repeated in C4 and a plain year value in E4. ---
$ws.Range("A4").Value = "hello ford mustang "
$ws.Range("A4").Font.Color = 255          # red (RGB 255,0,0)

$ws.Range("C4").Value = "new model"

$ws.Range("E4").Value = 2024

# --- Row 5: new status message ---
$ws.Range("A5").Value = "new message added"

# Leave the selection where the author left it when saving
$ws.Range("A5").Select()

# The saved file prints in portrait orientation
$ws.PageSetup.Orientation = 1
